$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Model"
$ws.Range("B1").Value = "Accuracy"
$ws.Range("C1").Value = "Precision"
$ws.Range("D1").Value = "Recall"
$ws.Range("E1").Value = "F1"

# Row 2 - Multi LR
$ws.Range("A2").Value = "Multi LR"
$ws.Range("B2").Value = 0.971
$ws.Range("C2").Value = 0.97
$ws.Range("D2").Value = 0.994
$ws.Range("E2").Value = 0.982

# Row 3 - Multi RF
$ws.Range("A3").Value = "Multi RF"
$ws.Range("B3").Value = 0.948
$ws.Range("C3").Value = 0.964
$ws.Range("D3").Value = 0.97
$ws.Range("E3").Value = 0.967

# Row 4 - Single LR
$ws.Range("A4").Value = "Single LR"
$ws.Range("B4").Value = 0.929
$ws.Range("C4").Value = 0.946
$ws.Range("D4").Value = 0.963
$ws.Range("E4").Value = 0.955

# Row 5 - Single RF
$ws.Range("A5").Value = "Single RF"
$ws.Range("B5").Value = 0.914
$ws.Range("C5").Value = 0.9399999999999999
$ws.Range("D5").Value = 0.951
$ws.Range("E5").Value = 0.945
